$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data ---
# B2 email changes to the new forgot-password test account
$ws.Range("B2").Value = "lakheramohini98@gmail.com"

# --- Row 1 headers (D1:F1 new, bold "header" style matching A1:C1) ---
$ws.Range("D1").Value = "pass1"
$ws.Range("E1").Value = "pass2"
$ws.Range("F1").Value = "password"
$ws.Range("D1:F1").Font.Bold = $true

$ws.Range("D2").Value = "Lakhera@123"
$ws.Range("E2").Value = "Lakhera@1234"
$ws.Range("F2").Value = "Lakhera@1998"

# --- Hyperlinks for the new password columns (matches A2/B2 pattern) ---
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Lakhera@123")
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Lakhera@1234")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:Lakhera@1998")

# Re-apply the shared "Hyperlink" look so the new cells reuse the same style as A2/B2
$ws.Range("D2:F2").Style = "Hyperlink"

# --- Active selection moves to B1 ---
$ws.Range("B1").Select()
